$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1917808219178082
$ws.Range("C2").Value = 0.5547945205479452
$ws.Range("J2").Value = 0.01712328767123288
$ws.Range("P2").Value = 0.1232876712328767
$ws.Range("S2").Value = 0.113013698630137
$ws.Range("B3").Value = 0.01829268292682927
$ws.Range("C3").Value = 0.01829268292682927
$ws.Range("J3").Value = 0.01829268292682927
$ws.Range("P3").Value = 0.7378048780487805
$ws.Range("S3").Value = 0.2073170731707317
$ws.Range("J4").Value = 0.07894736842105263
$ws.Range("P4").Value = 0.6578947368421053
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.05579399141630902
$ws.Range("F6").Value = 0.09871244635193133
$ws.Range("J6").Value = 0.240343347639485
$ws.Range("O6").Value = 0.04721030042918455
$ws.Range("Q6").Value = 0.1330472103004292
$ws.Range("R6").Value = 0.06437768240343347
$ws.Range("S6").Value = 0.3605150214592275
$ws.Range("B7").Value = 0.1059907834101382
$ws.Range("D7").Value = 0.0184331797235023
$ws.Range("F7").Value = 0.07373271889400922
$ws.Range("J7").Value = 0.1797235023041475
$ws.Range("O7").Value = 0.02304147465437788
$ws.Range("Q7").Value = 0.1474654377880184
$ws.Range("R7").Value = 0.08755760368663594
$ws.Range("S7").Value = 0.3640552995391705
$ws.Range("B8").Value = 0.05427974947807934
$ws.Range("D8").Value = 0.02922755741127349
$ws.Range("E8").Value = 0.004175365344467641
$ws.Range("F8").Value = 0.08559498956158663
$ws.Range("J8").Value = 0.1169102296450939
$ws.Range("O8").Value = 0.02296450939457203
$ws.Range("Q8").Value = 0.1670146137787056
$ws.Range("R8").Value = 0.08350730688935282
$ws.Range("S8").Value = 0.4363256784968685
$ws.Range("B9").Value = 0.09905660377358491
$ws.Range("D9").Value = 0.009433962264150943
$ws.Range("F9").Value = 0.07547169811320754
$ws.Range("J9").Value = 0.1650943396226415
$ws.Range("O9").Value = 0.02830188679245283
$ws.Range("Q9").Value = 0.1367924528301887
$ws.Range("R9").Value = 0.1462264150943396
$ws.Range("S9").Value = 0.3396226415094339
$ws.Range("B10").Value = 0.1083969465648855
$ws.Range("D10").Value = 0.01526717557251908
$ws.Range("F10").Value = 0.05801526717557252
$ws.Range("J10").Value = 0.1343511450381679
$ws.Range("O10").Value = 0.01450381679389313
$ws.Range("Q10").Value = 0.1877862595419847
$ws.Range("R10").Value = 0.09923664122137404
$ws.Range("S10").Value = 0.3824427480916031
$ws.Range("G11").Value = 0.1223880597014925
$ws.Range("J11").Value = 0.1104477611940299
$ws.Range("K11").Value = 0.173134328358209
$ws.Range("L11").Value = 0.5880597014925373
$ws.Range("S11").Value = 0.005970149253731343
$ws.Range("G12").Value = 0.76
$ws.Range("J12").Value = 0.195
$ws.Range("K12").Value = 0.005
$ws.Range("L12").Value = 0.015
$ws.Range("S12").Value = 0.025
$ws.Range("F13").Value = 0.02
$ws.Range("G13").Value = 0.5600000000000001
$ws.Range("J13").Value = 0.38
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.02242152466367713
$ws.Range("H15").Value = 0.2017937219730942
$ws.Range("I15").Value = 0.07623318385650224
$ws.Range("J15").Value = 0.2376681614349776
$ws.Range("K15").Value = 0.08968609865470852
$ws.Range("M15").Value = 0.008968609865470852
$ws.Range("O15").Value = 0.09865470852017937
$ws.Range("S15").Value = 0.2645739910313901
$ws.Range("F16").Value = 0.01123595505617977
$ws.Range("H16").Value = 0.1573033707865168
$ws.Range("I16").Value = 0.101123595505618
$ws.Range("J16").Value = 0.4213483146067415
$ws.Range("K16").Value = 0.1067415730337079
$ws.Range("M16").Value = 0.02247191011235955
$ws.Range("O16").Value = 0.02247191011235955
$ws.Range("S16").Value = 0.1573033707865168
$ws.Range("F17").Value = 0.01909307875894988
$ws.Range("H17").Value = 0.1766109785202864
$ws.Range("I17").Value = 0.09307875894988067
$ws.Range("J17").Value = 0.3985680190930788
$ws.Range("K17").Value = 0.1145584725536993
$ws.Range("M17").Value = 0.009546539379474941
$ws.Range("O17").Value = 0.05966587112171837
$ws.Range("S17").Value = 0.1288782816229117
$ws.Range("F18").Value = 0.02978723404255319
$ws.Range("H18").Value = 0.2127659574468085
$ws.Range("I18").Value = 0.1148936170212766
$ws.Range("J18").Value = 0.3787234042553191
$ws.Range("K18").Value = 0.1021276595744681
$ws.Range("M18").Value = 0.01276595744680851
$ws.Range("O18").Value = 0.02978723404255319
$ws.Range("S18").Value = 0.1191489361702128
$ws.Range("F19").Value = 0.013215859030837
$ws.Range("H19").Value = 0.2085168869309839
$ws.Range("I19").Value = 0.07929515418502203
$ws.Range("J19").Value = 0.3568281938325991
$ws.Range("K19").Value = 0.1196769456681351
$ws.Range("M19").Value = 0.02936857562408223
$ws.Range("N19").Value = 0.002202643171806168
$ws.Range("O19").Value = 0.06240822320117474
$ws.Range("S19").Value = 0.1284875183553598
